$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column B (1989 data) - this also removes the formatting stub cells B1/B2/B3
$ws.Columns.Item(2).Delete()
# Delete column B again (now holds 2002 data)
$ws.Columns.Item(2).Delete()
# Delete the blank row 3
$ws.Rows.Item(3).Delete()

# Clear the explanatory note text (and its formatting) in A2
$ws.Range("A2").Clear()

# Set the uniform custom row height used after the cleanup
$ws.Rows.Item(1).RowHeight = 20.1
$ws.Rows.Item(2).RowHeight = 20.1
$ws.Rows.Item(3).RowHeight = 20.1
$ws.Rows.Item(4).RowHeight = 20.1
$ws.Rows.Item(5).RowHeight = 20.1
